$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1707317073170732
$ws.Range("C2").Value = 0.6097560975609756
$ws.Range("J2").Value = 0.02439024390243903
$ws.Range("P2").Value = 0.1219512195121951
$ws.Range("S2").Value = 0.07317073170731707
$ws.Range("B3").Value = 0.005649717514124294
$ws.Range("C3").Value = 0.02259887005649718
$ws.Range("J3").Value = 0.02824858757062147
$ws.Range("P3").Value = 0.7457627118644068
$ws.Range("S3").Value = 0.1977401129943503
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.7045454545454546
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.06930693069306931
$ws.Range("D6").Value = 0.009900990099009901
$ws.Range("F6").Value = 0.1138613861386139
$ws.Range("J6").Value = 0.2277227722772277
$ws.Range("O6").Value = 0.0198019801980198
$ws.Range("Q6").Value = 0.1237623762376238
$ws.Range("R6").Value = 0.06435643564356436
$ws.Range("S6").Value = 0.3712871287128713
$ws.Range("B7").Value = 0.1136363636363636
$ws.Range("D7").Value = 0.01136363636363636
$ws.Range("F7").Value = 0.03409090909090909
$ws.Range("J7").Value = 0.1818181818181818
$ws.Range("O7").Value = 0.01704545454545454
$ws.Range("Q7").Value = 0.2215909090909091
$ws.Range("R7").Value = 0.08522727272727272
$ws.Range("S7").Value = 0.3352272727272727
$ws.Range("B8").Value = 0.09788359788359788
$ws.Range("D8").Value = 0.005291005291005291
$ws.Range("F8").Value = 0.06349206349206349
$ws.Range("J8").Value = 0.1005291005291005
$ws.Range("O8").Value = 0.02116402116402116
$ws.Range("Q8").Value = 0.1931216931216931
$ws.Range("R8").Value = 0.126984126984127
$ws.Range("S8").Value = 0.3915343915343915
$ws.Range("B9").Value = 0.1074766355140187
$ws.Range("D9").Value = 0.03271028037383177
$ws.Range("E9").Value = 0.004672897196261682
$ws.Range("F9").Value = 0.0514018691588785
$ws.Range("J9").Value = 0.1588785046728972
$ws.Range("O9").Value = 0.009345794392523364
$ws.Range("Q9").Value = 0.1635514018691589
$ws.Range("R9").Value = 0.05607476635514019
$ws.Range("S9").Value = 0.4158878504672897
$ws.Range("B10").Value = 0.1091908876669285
$ws.Range("D10").Value = 0.02592301649646504
$ws.Range("E10").Value = 0.001571091908876669
$ws.Range("F10").Value = 0.06048703849175176
$ws.Range("J10").Value = 0.1500392772977219
$ws.Range("O10").Value = 0.01649646504320503
$ws.Range("Q10").Value = 0.2073841319717203
$ws.Range("R10").Value = 0.08562450903377848
$ws.Range("S10").Value = 0.3432835820895522
$ws.Range("G11").Value = 0.149812734082397
$ws.Range("J11").Value = 0.0898876404494382
$ws.Range("K11").Value = 0.2097378277153558
$ws.Range("L11").Value = 0.5318352059925093
$ws.Range("S11").Value = 0.01872659176029963
$ws.Range("G12").Value = 0.7414965986394558
$ws.Range("J12").Value = 0.1836734693877551
$ws.Range("K12").Value = 0.006802721088435374
$ws.Range("L12").Value = 0.02040816326530612
$ws.Range("S12").Value = 0.04761904761904762
$ws.Range("F13").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.6279069767441861
$ws.Range("J13").Value = 0.3023255813953488
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.01345291479820628
$ws.Range("H15").Value = 0.09865470852017937
$ws.Range("I15").Value = 0.06278026905829596
$ws.Range("J15").Value = 0.3946188340807175
$ws.Range("K15").Value = 0.08071748878923767
$ws.Range("M15").Value = 0.01345291479820628
$ws.Range("O15").Value = 0.1076233183856502
$ws.Range("S15").Value = 0.2286995515695067
$ws.Range("F16").Value = 0.02116402116402116
$ws.Range("H16").Value = 0.164021164021164
$ws.Range("I16").Value = 0.07407407407407407
$ws.Range("J16").Value = 0.4232804232804233
$ws.Range("K16").Value = 0.1216931216931217
$ws.Range("M16").Value = 0.02116402116402116
$ws.Range("N16").Value = 0.005291005291005291
$ws.Range("O16").Value = 0.1005291005291005
$ws.Range("S16").Value = 0.06878306878306878
$ws.Range("F17").Value = 0.0162037037037037
$ws.Range("H17").Value = 0.1643518518518519
$ws.Range("I17").Value = 0.1180555555555556
$ws.Range("J17").Value = 0.4074074074074074
$ws.Range("K17").Value = 0.1064814814814815
$ws.Range("M17").Value = 0.0162037037037037
$ws.Range("O17").Value = 0.06018518518518518
$ws.Range("S17").Value = 0.1111111111111111
$ws.Range("F18").Value = 0.03061224489795918
$ws.Range("H18").Value = 0.173469387755102
$ws.Range("I18").Value = 0.07653061224489796
$ws.Range("J18").Value = 0.4438775510204082
$ws.Range("K18").Value = 0.05612244897959184
$ws.Range("M18").Value = 0.02040816326530612
$ws.Range("O18").Value = 0.05612244897959184
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.01533219761499148
$ws.Range("H19").Value = 0.182282793867121
$ws.Range("I19").Value = 0.100511073253833
$ws.Range("J19").Value = 0.3824531516183987
$ws.Range("K19").Value = 0.09369676320272573
$ws.Range("M19").Value = 0.02129471890971039
$ws.Range("O19").Value = 0.06814310051107325
$ws.Range("S19").Value = 0.1362862010221465
